$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# "Generate Report for Handback" - fill in the handback-received columns
# for f0b0f9f7-cc12-4ad3-8892-dcea71ee903c (row 7) on both the zh-cn and
# de-de sheets: a handback file/url was picked up, but it is not the
# latest version, so status column is not changed but Latest Target
# File / Latest Handback File / Latest Handback DateTime / Error Detail
# get populated.
# -------------------------------------------------------------------------

$mdName        = "f0b0f9f7-cc12-4ad3-8892-dcea71ee903c.md"
$mdUrl         = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/962fb2ddb794ec09a66140012337dc8673f6ec8b/e2e/f0b0f9f7-cc12-4ad3-8892-dcea71ee903c.md"
$handbackDate  = "2016-09-04 08:59:09"
$errorDetail   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/927582fc730f6b2c390b866662f4beedbc7a687a/e2e/f0b0f9f7-cc12-4ad3-8892-dcea71ee903c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/962fb2ddb794ec09a66140012337dc8673f6ec8b/e2e/f0b0f9f7-cc12-4ad3-8892-dcea71ee903c.md."

# ---------------------- zh-cn sheet (row 7) -----------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsZh.Range("J7").Value = "f0b0f9f7-cc12-4ad3-8892-dcea71ee903c.86617a4ad85cd40b8d5e30646d8792284efd94d8.zh-cn.xlf"
$wsZh.Range("K7").Value = $handbackDate
$wsZh.Range("P7").Value = $errorDetail

# ---------------------- de-de sheet (row 7) -----------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsDe.Range("J7").Value = "f0b0f9f7-cc12-4ad3-8892-dcea71ee903c.86617a4ad85cd40b8d5e30646d8792284efd94d8.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-04 08:59:16"
$wsDe.Range("P7").Value = $errorDetail
